{"js": "// Fill in the lab cover-page table: Drive Full Name, Partner Full Name, Student ID.\nconst table = context.document.body.tables.getFirst();\ntable.load(\"rowCount\");\nawait context.sync();\n\n// Row 0 -> \"Drive Full Name\" value cell\nconst driveCell = table.getCell(0, 1);\ndriveCell.body.insertText(\"Theresa DeJacimo\", \"End\");\n\n// Row 1 -> \"Partner Full Name\" value cell\nconst partnerCell = table.getCell(1, 1);\npartnerCell.body.insertText(\"Cameron Combariza \", \"End\");\n\n// Row 2 -> \"Student ID\" value cell\nconst idCell = table.getCell(2, 1);\nidCell.body.insertText(\"1910844\", \"End\");\n\nawait context.sync();\n", "ps1": "# Fill in the lab cover-page table: Drive Full Name, Partner Full Name, Student ID.\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Row 1 -> \"Drive Full Name\" value cell\n$driveRange = $table.Cell(1, 2).Range\n$driveRange.Collapse(0)\n$driveRange.InsertAfter(\"Theresa DeJacimo\")\n\n# Row 2 -> \"Partner Full Name\" value cell\n$partnerRange = $table.Cell(2, 2).Range\n$partnerRange.Collapse(0)\n$partnerRange.InsertAfter(\"Cameron Combariza \")\n\n# Row 3 -> \"Student ID\" value cell\n$idRange = $table.Cell(3, 2).Range\n$idRange.Collapse(0)\n$idRange.InsertAfter(\"1910844\")\n"}
